# Add "Wins", "Losses", "Ties" columns (season record) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new column titles in AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style used by the other header cells (e.g. AC1) onto the new headers
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

# Data rows 2-53: season record values (Wins=66, Losses=96, Ties=0)
$lastRow = 53
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 66
    $ws.Cells.Item($r, 31).Value = 96
    $ws.Cells.Item($r, 32).Value = 0
}
